# Update the "Metadata" sheet (Property/Value pairs in column A/B).
$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-primary-physician"  # URL
$meta.Range("B3").Value = "8.0.0"                                                                              # Version
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"                                                          # Date
$meta.Range("B9").Value = "LinuxForHealth Team"                                                                # Publisher

# Update the "Elements" sheet: clear the stale root-level FHIR constraint text
# (row 2 = "Extension" element, column AI = "Constraint(s)").
$elements = $wb.Worksheets("Elements")
$elements.Range("AI2").Value = ""

# Extension.url's fixed value (row 5, column Q) mirrors the canonical URL shown
# on the Metadata sheet, so it must be kept in sync with the new URL as well.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-primary-physician"
